$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Platform Coverage")
$ws2 = $wb.Worksheets.Item("MarketShare")

# --- Platform Coverage sheet -------------------------------------------------
# Row 2 (All / Treatment / Campaign / MDA): extend the 0.6 coverage values to
# every other year from 2026 through 2040 (columns P,R,T,V,X,Z,AB,AD).
$ws1.Range("P2").Value = 0.6
$ws1.Range("R2").Value = 0.6
$ws1.Range("T2").Value = 0.6
$ws1.Range("V2").Value = 0.6
$ws1.Range("X2").Value = 0.6
$ws1.Range("Z2").Value = 0.6
$ws1.Range("AB2").Value = 0.6
$ws1.Range("AD2").Value = 0.6

# Old row 3 (All / Treatment / Campaign / MDA, age 0-15 w/ the 0.7 overrides)
# is removed entirely; rows 4-8 shift up to become rows 3-7.
$ws1.Rows.Item(3).Delete()

# --- MarketShare sheet -------------------------------------------------------
# New Product A's market share (row 2) no longer ramps up from 2026 (L2:Z2);
# Old Product B (SOC) (row 3) instead keeps its full 100% share for those
# years too, extending its existing D3:K3 = 1 run through to Z3.
$ws2.Range("L3:Z3").Value = 1
$ws2.Range("L2:Z2").ClearContents()

# --- Restore view/selection state -------------------------------------------
$ws2.Activate()
$ws2.Range("Z3").Select()

$ws1.Activate()
$ws1.Range("G6").Select()
